# Insert one new data row just above the existing "2026/12/29" block.
# Before: row 854 = 2026/12/29 (the first of several rows for that date),
# preceded by a single leftover row 853 = 2026/02/23 (only one time slot
# logged so far for that day). The edit adds the missing second time
# slot for 2026/02/23 (時刻=5) as a new row 854, pushing every following
# row down by one (854..895 -> 855..896) and extending the sheet by one
# row (A1:D895 -> A1:D896).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 854:895 down to 855:896, opening up a blank row 854.
$ws.Rows("854:854").Insert()

# Fill the new row. Column A holds a date written as plain text
# ("2026/02/23"), matching how every other date cell in this sheet is
# stored (t="inlineStr"/shared-string text, not a real date serial).
# Flip the cell to text format before assigning so the engine doesn't
# auto-convert the string into a date number, then restore the cell's
# style to Normal so no stray style index is left behind (the original
# cells in this column carry no explicit style).
$newRow = 854
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2026/02/23"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "月"
$ws.Cells.Item($newRow, 3).Value = 5
$ws.Cells.Item($newRow, 4).Value = 201
